# "Add tester to test server"
# The MasterServer_1 row's internal IP is repointed from the old
# 192.168.150.100 address to the local tester address 127.0.0.1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("property")

# Row 11 = MasterServer_1 -> IP column (E) becomes the local tester IP.
$ws.Range("E11").Value = "127.0.0.1"

# Minor re-layout that rode along with the save (header row got a touch
# shorter once the new value was in place).
$ws.Rows.Item(10).RowHeight = 28.35

# Leave the cursor where the author left it when they saved the file.
[void]$ws.Range("F30").Select()
